# Apply updated symbol list values (price, volume, hora) to the crypto sheet.
# All target columns (D=Price, E=Volume(1h), G=Hora) are stored as text in the
# worksheet, so we temporarily switch each cell to the "Text" number format
# before writing the value (otherwise Excel auto-converts numeric-looking or
# percentage-looking strings into real numbers), then restore the cell's
# original style so no stray formatting is introduced.
function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "308.60"
Set-TextValue $ws "E2" "-2.85%"
Set-TextValue $ws "G2" "1"
Set-TextValue $ws "D3" "37.51"
Set-TextValue $ws "E3" "-4.21%"
Set-TextValue $ws "G3" "1"
Set-TextValue $ws "D4" "5.068"
Set-TextValue $ws "E4" "-1.49%"
Set-TextValue $ws "G4" "1"
Set-TextValue $ws "D5" "0.07762"
Set-TextValue $ws "E5" "-5.22%"
Set-TextValue $ws "G5" "1"
Set-TextValue $ws "D6" "4.344"
Set-TextValue $ws "E6" "0.02%"
Set-TextValue $ws "G6" "1"
Set-TextValue $ws "D7" "1.902"
Set-TextValue $ws "E7" "-5.52%"
Set-TextValue $ws "G7" "1"
Set-TextValue $ws "D8" "8.180"
Set-TextValue $ws "E8" "-2.18%"
Set-TextValue $ws "G8" "1"
Set-TextValue $ws "D9" "3.094"
Set-TextValue $ws "E9" "-6.37%"
Set-TextValue $ws "G9" "1"
Set-TextValue $ws "D10" "0.9216"
Set-TextValue $ws "E10" "-1.86%"
Set-TextValue $ws "G10" "1"
Set-TextValue $ws "D11" "0.1257"
Set-TextValue $ws "E11" "-3.50%"
Set-TextValue $ws "G11" "1"
Set-TextValue $ws "D12" "0.1863"
Set-TextValue $ws "E12" "-6.32%"
Set-TextValue $ws "G12" "1"
Set-TextValue $ws "D13" "0.08748"
Set-TextValue $ws "E13" "-3.73%"
Set-TextValue $ws "G13" "1"
Set-TextValue $ws "D14" "0.03384"
Set-TextValue $ws "E14" "-3.43%"
Set-TextValue $ws "G14" "1"
Set-TextValue $ws "D15" "0.09683"
Set-TextValue $ws "E15" "-0.81%"
Set-TextValue $ws "G15" "1"
Set-TextValue $ws "D16" "0.001371"
Set-TextValue $ws "E16" "-2.69%"
Set-TextValue $ws "G16" "1"
Set-TextValue $ws "D17" "0.005800"
Set-TextValue $ws "E17" "-3.39%"
Set-TextValue $ws "G17" "1"
Set-TextValue $ws "D18" "3.600"
Set-TextValue $ws "E18" "-2.50%"
Set-TextValue $ws "G18" "1"
Set-TextValue $ws "D19" "0.3410"
Set-TextValue $ws "E19" "-2.30%"
Set-TextValue $ws "G19" "1"
Set-TextValue $ws "B20" "MCDex"
Set-TextValue $ws "C20" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws "D20" "5.027"
Set-TextValue $ws "E20" "1.29%"
Set-TextValue $ws "G20" "1"
Set-TextValue $ws "B21" "ProBitToken"
Set-TextValue $ws "C21" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws "D21" "0.1266"
Set-TextValue $ws "E21" "-3.85%"
Set-TextValue $ws "G21" "1"
Set-TextValue $ws "D22" "0.2570"
Set-TextValue $ws "E22" "2.49%"
Set-TextValue $ws "G22" "1"
Set-TextValue $ws "D23" "0.02107"
Set-TextValue $ws "E23" "5,596.40%"
Set-TextValue $ws "G23" "1"
Set-TextValue $ws "D24" "0.04311"
Set-TextValue $ws "E24" "-0.69%"
Set-TextValue $ws "G24" "1"
Set-TextValue $ws "D25" "0.001212"
Set-TextValue $ws "E25" "-2.47%"
Set-TextValue $ws "G25" "1"
Set-TextValue $ws "D26" "0.004235"
Set-TextValue $ws "E26" "-11.12%"
Set-TextValue $ws "G26" "1"
Set-TextValue $ws "D27" "0.0001351"
Set-TextValue $ws "E27" "3.81%"
Set-TextValue $ws "G27" "1"
Set-TextValue $ws "G28" "1"
Set-TextValue $ws "G29" "1"
Set-TextValue $ws "G30" "1"
Set-TextValue $ws "G31" "1"
Set-TextValue $ws "G32" "1"
Set-TextValue $ws "G33" "1"
Set-TextValue $ws "G34" "1"
Set-TextValue $ws "G35" "1"
Set-TextValue $ws "G36" "1"
Set-TextValue $ws "G37" "1"
Set-TextValue $ws "G38" "1"
Set-TextValue $ws "D39" "0.02140"
Set-TextValue $ws "E39" "-5.28%"
Set-TextValue $ws "G39" "1"
Set-TextValue $ws "D40" "0.04922"
Set-TextValue $ws "E40" "-4.98%"
Set-TextValue $ws "G40" "1"
Set-TextValue $ws "D41" "0.007936"
Set-TextValue $ws "E41" "2.28%"
Set-TextValue $ws "G41" "1"
Set-TextValue $ws "D42" "0.01002"
Set-TextValue $ws "E42" "-4.63%"
Set-TextValue $ws "G42" "1"
Set-TextValue $ws "D43" "0.1339"
Set-TextValue $ws "E43" "-4.48%"
Set-TextValue $ws "G43" "1"
Set-TextValue $ws "D44" "0.001993"
Set-TextValue $ws "E44" "-2.57%"
Set-TextValue $ws "G44" "1"
Set-TextValue $ws "D45" "0.009570"
Set-TextValue $ws "E45" "5.85%"
Set-TextValue $ws "G45" "1"
Set-TextValue $ws "D46" "0.00006441"
Set-TextValue $ws "E46" "-7.04%"
Set-TextValue $ws "G46" "1"
Set-TextValue $ws "D47" "0.00000000751"
Set-TextValue $ws "E47" "-0.07%"
Set-TextValue $ws "G47" "1"
Set-TextValue $ws "D48" "0.003367"
Set-TextValue $ws "E48" "16.66%"
Set-TextValue $ws "G48" "1"
Set-TextValue $ws "D49" "0.001691"
Set-TextValue $ws "E49" "-0.12%"
Set-TextValue $ws "G49" "1"
Set-TextValue $ws "D50" "0.00002102"
Set-TextValue $ws "E50" "-0.07%"
Set-TextValue $ws "G50" "1"
Set-TextValue $ws "D51" "0.0002002"
Set-TextValue $ws "E51" "-0.07%"
Set-TextValue $ws "G51" "1"
